$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 8: update Inscritos, Pagos, Isenções deferidas, Inscrições homologadas
$ws.Range("E8").Value = 44
$ws.Range("F8").Value = 15
$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 21

# Row 14: update Inscritos, Isenções deferidas, Inscrições homologadas (Pagos unchanged)
$ws.Range("E14").Value = 36
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 18
